$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Fill in the previously-empty "Neural Network" row (row 9) with its
# timing / performance / confusion-matrix results.
$ws.Range("B9").Value = "1 minuto"
$ws.Range("C9").Value = "0.73 secondi"
$ws.Range("D9").Value = "0.835"
$ws.Range("E9").Value = "0.6688"
$ws.Range("F9").Value = "0.5205"
$ws.Range("G9").Value = "0.5854"
$ws.Range("H9").Value = "TP=3942, FP=1952, FN=3631, TN=24305"

# Update the view: scroll so column D is the left-most visible column,
# and move the active selection to H9.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H9").Select()
